# Add newly collected landscaping data rows (359-365) to the bottom of the
# existing data table on sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 359
$endRow = 365

# Carry the formatting (date number format on column A, general elsewhere)
# from the last existing row down onto the new rows before filling values,
# so no new style entries are created.
$ws.Range("A358:T358").Copy() | Out-Null
$ws.Range("A$startRow`:T$endRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data for the new rows, in column order A..T
# A: Date, B: Plant_Type, C: Plant_Size, D: Low, E: High, F: Temp_Diff (formula),
# G: Rain, H: Growth, I: Pruned, J: Quadrant, K: Shade, L: UV, M: Humidity,
# N: Dew_Point, O: Pressure, P: Wind_Gust, Q: Cloud_Cover, R: Visibility,
# S: AQI, T: Pollen

$newRows = @(
    @(45838, "Flowering",    "Large",  72, 84, 0.29, 0.3,  "Yes", 2, "Bright",  7, 0.74, 74, 29.98, 12, 0.62, 9.3, 53, 0),
    @(45838, "Nonflowering", "Medium", 72, 84, 0.29, 0.35, "Yes", 3, "Bright",  7, 0.74, 74, 29.98, 12, 0.62, 9.3, 53, 0),
    @(45838, "Nonflowering", "Small",  72, 84, 0.29, 0.35, "Yes", 3, "Neutral", 7, 0.74, 74, 29.98, 12, 0.62, 9.3, 53, 0),
    @(45838, "Nonflowering", "Medium", 72, 84, 0.29, 0.5,  "Yes", 3, "Neutral", 7, 0.74, 74, 29.98, 12, 0.62, 9.3, 53, 0),
    @(45838, "Nonflowering", "Medium", 72, 84, 0.29, 0.6,  "Yes", 3, "Bright",  7, 0.74, 74, 29.98, 12, 0.62, 9.3, 53, 0),
    @(45838, "Nonflowering", "Large",  72, 84, 0.29, 0.7,  "Yes", 4, "Neutral", 7, 0.74, 74, 29.98, 12, 0.62, 9.3, 53, 0),
    @(45838, "Tree",         "Medium", 72, 84, 0.29, 1.3,  "Yes", 1, "Bright",  7, 0.74, 74, 29.98, 12, 0.62, 9.3, 53, 0)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Formula = "=ABS(D$r-E$r)"
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[8]
    $ws.Cells.Item($r, 11).Value = $row[9]
    $ws.Cells.Item($r, 12).Value = $row[10]
    $ws.Cells.Item($r, 13).Value = $row[11]
    $ws.Cells.Item($r, 14).Value = $row[12]
    $ws.Cells.Item($r, 15).Value = $row[13]
    $ws.Cells.Item($r, 16).Value = $row[14]
    $ws.Cells.Item($r, 17).Value = $row[15]
    $ws.Cells.Item($r, 18).Value = $row[16]
    $ws.Cells.Item($r, 19).Value = $row[17]
    $ws.Cells.Item($r, 20).Value = $row[18]
}

# Update the view to reflect where the user was working after the edit
$excel.ActiveWindow.ScrollRow = 343
$ws.Range("G366").Select()
